# Update TPM-derived NATMI ligand-receptor metrics with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.02328126719340038
$ws.Range("J2").Value = 0.02328126719340038
$ws.Range("M2").Value = 1.090710333333333
$ws.Range("N2").Value = 3.272131
$ws.Range("O2").Value = 0.0488470045579656
$ws.Range("P2").Value = 0.0488470045579656
$ws.Range("Q2").Value = 0.022999808799
$ws.Range("R2").Value = 0.206998279191
$ws.Range("S2").Value = 0.001137220164711243
$ws.Range("T2").Value = 0.001137220164711244

# Row 3
$ws.Range("I3").Value = 0.02328126719340038
$ws.Range("J3").Value = 0.02328126719340038
$ws.Range("O3").Value = 0.7616320856558244
$ws.Range("P3").Value = 0.7616320856558244
$ws.Range("S3").Value = 0.01773176008922005
$ws.Range("T3").Value = 0.01773176008922006

# Row 4
$ws.Range("I4").Value = 0.02328126719340038
$ws.Range("J4").Value = 0.02328126719340038
$ws.Range("M4").Value = 0.740281
$ws.Range("N4").Value = 2.220843
$ws.Range("O4").Value = 0.03315317392351528
$ws.Range("P4").Value = 0.03315317392351528
$ws.Range("Q4").Value = 0.015610305447
$ws.Range("R4").Value = 0.140492749023
$ws.Range("S4").Value = 0.0007718479004226332
$ws.Range("T4").Value = 0.0007718479004226334

# Row 5
$ws.Range("I5").Value = 0.02328126719340038
$ws.Range("J5").Value = 0.02328126719340038
$ws.Range("M5").Value = 2.784013333333333
$ws.Range("N5").Value = 8.352039999999999
$ws.Range("O5").Value = 0.1246808688124989
$ws.Range("P5").Value = 0.1246808688124989
$ws.Range("Q5").Value = 0.05870648915999999
$ws.Range("R5").Value = 0.5283584024399999
$ws.Range("S5").Value = 0.002902728620729087
$ws.Range("T5").Value = 0.002902728620729088

# Row 6
$ws.Range("I6").Value = 0.02328126719340038
$ws.Range("J6").Value = 0.02328126719340038
$ws.Range("M6").Value = 0.2710316666666667
$ws.Range("N6").Value = 0.813095
$ws.Range("O6").Value = 0.01213803945228936
$ws.Range("P6").Value = 0.01213803945228936
$ws.Range("Q6").Value = 0.005715244755
$ws.Range("R6").Value = 0.051437202795
$ws.Range("S6").Value = 0.0002825889396927838
$ws.Range("T6").Value = 0.0002825889396927838

# Row 7
$ws.Range("I7").Value = 0.02328126719340038
$ws.Range("J7").Value = 0.02328126719340038
$ws.Range("M7").Value = 0.436508
$ws.Range("N7").Value = 1.309524
$ws.Range("O7").Value = 0.01954882759790648
$ws.Range("P7").Value = 0.01954882759790648
$ws.Range("Q7").Value = 0.009204644195999999
$ws.Range("R7").Value = 0.082841797764
$ws.Range("S7").Value = 0.0004551214786245801
$ws.Range("T7").Value = 0.0004551214786245802

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8846626666666667
$ws.Range("H8").Value = 2.653988
$ws.Range("I8").Value = 0.9767187328065996
$ws.Range("J8").Value = 0.9767187328065997
$ws.Range("M8").Value = 1.090710333333333
$ws.Range("N8").Value = 3.272131
$ws.Range("O8").Value = 0.0488470045579656
$ws.Range("P8").Value = 0.0488470045579656
$ws.Range("Q8").Value = 0.9649107120475557
$ws.Range("R8").Value = 8.684196408428001
$ws.Range("S8").Value = 0.04770978439325435
$ws.Range("T8").Value = 0.04770978439325436

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8846626666666667
$ws.Range("H9").Value = 2.653988
$ws.Range("I9").Value = 0.9767187328065996
$ws.Range("J9").Value = 0.9767187328065997
$ws.Range("O9").Value = 0.7616320856558244
$ws.Range("P9").Value = 0.7616320856558244
$ws.Range("Q9").Value = 15.04507727216578
$ws.Range("R9").Value = 135.405695449492
$ws.Range("S9").Value = 0.7439003255666043
$ws.Range("T9").Value = 0.7439003255666045

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8846626666666667
$ws.Range("H10").Value = 2.653988
$ws.Range("I10").Value = 0.9767187328065996
$ws.Range("J10").Value = 0.9767187328065997
$ws.Range("M10").Value = 0.740281
$ws.Range("N10").Value = 2.220843
$ws.Range("O10").Value = 0.03315317392351528
$ws.Range("P10").Value = 0.03315317392351528
$ws.Range("Q10").Value = 0.6548989635426666
$ws.Range("R10").Value = 5.894090671883999
$ws.Range("S10").Value = 0.03238132602309265
$ws.Range("T10").Value = 0.03238132602309265

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.8846626666666667
$ws.Range("H11").Value = 2.653988
$ws.Range("I11").Value = 0.9767187328065996
$ws.Range("J11").Value = 0.9767187328065997
$ws.Range("M11").Value = 2.784013333333333
$ws.Range("N11").Value = 8.352039999999999
$ws.Range("O11").Value = 0.1246808688124989
$ws.Range("P11").Value = 0.1246808688124989
$ws.Range("Q11").Value = 2.462912659502222
$ws.Range("R11").Value = 22.16621393552
$ws.Range("S11").Value = 0.1217781401917698
$ws.Range("T11").Value = 0.1217781401917698

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.8846626666666667
$ws.Range("H12").Value = 2.653988
$ws.Range("I12").Value = 0.9767187328065996
$ws.Range("J12").Value = 0.9767187328065997
$ws.Range("M12").Value = 0.2710316666666667
$ws.Range("N12").Value = 0.813095
$ws.Range("O12").Value = 0.01213803945228936
$ws.Range("P12").Value = 0.01213803945228936
$ws.Range("Q12").Value = 0.2397715969844444
$ws.Range("R12").Value = 2.15794437286
$ws.Range("S12").Value = 0.01185545051259657
$ws.Range("T12").Value = 0.01185545051259658

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.8846626666666667
$ws.Range("H13").Value = 2.653988
$ws.Range("I13").Value = 0.9767187328065996
$ws.Range("J13").Value = 0.9767187328065997
$ws.Range("M13").Value = 0.436508
$ws.Range("N13").Value = 1.309524
$ws.Range("O13").Value = 0.01954882759790648
$ws.Range("P13").Value = 0.01954882759790648
$ws.Range("Q13").Value = 0.3861623313013333
$ws.Range("R13").Value = 3.475460981712
$ws.Range("S13").Value = 0.0190937061192819
$ws.Range("T13").Value = 0.01909370611928191
